$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version 5 Created picture dataset for siamese
# Re-pair the near-duplicate "siamese" rows (File1Name groups that share the
# same HashFile1) with their correct File2Name/HashFile2/Match partners.

$ws.Range("D2").Value = '9ada1c93f3ed500fcdcee67415efe684'
$ws.Range("D3").Value = '9ada1c93f3ed500fcdcee67415efe684'
$ws.Range("D4").Value = '7cc9db896fc77886d4382ee98cc217db'
$ws.Range("D5").Value = '7cc9db896fc77886d4382ee98cc217db'
$ws.Range("B36").Value = 'fio3.sh'
$ws.Range("D36").Value = 'f995d8a0bc8d1f2528220691f5156c37'
$ws.Range("B37").Value = 'fio4.sh'
$ws.Range("D37").Value = 'fcc11badd30b5102a26311a41b951a3c'
$ws.Range("B38").Value = 'fio7.sh'
$ws.Range("D38").Value = '93d3ec88c2682a6b17d0a0b845c0b772'
$ws.Range("B39").Value = 'fio6.sh'
$ws.Range("D39").Value = '6376420d566d841ce003e14021f5c792'
$ws.Range("B40").Value = 'fio2.sh'
$ws.Range("D40").Value = '512b9068b4cd702237ab8242caceab38'
$ws.Range("B41").Value = 'fio5.sh'
$ws.Range("D41").Value = '9040171748b12ce1ef6a8536ac982b3e'
$ws.Range("B42").Value = 'fio9.sh'
$ws.Range("D42").Value = '61183c5e98db4661b70386d3bfc25ad1'
$ws.Range("B43").Value = 'fio8.sh'
$ws.Range("D43").Value = 'd436cb0e1a476937a5a1957ea19a530a'
$ws.Range("B44").Value = 'fio15.sh'
$ws.Range("D44").Value = '0e7fefd4bf12e3b10921488a1733c871'
$ws.Range("B45").Value = 'fio16.sh'
$ws.Range("D45").Value = '172fa5348233549194189081d49f9f6d'
$ws.Range("B46").Value = 'fio17.sh'
$ws.Range("D46").Value = '5017eef8adf42f9837bbd372d66e547d'
$ws.Range("B47").Value = 'fio11.sh'
$ws.Range("D47").Value = 'e22aa548fea442a35b145881113905b5'
$ws.Range("B49").Value = 'fio14.sh'
$ws.Range("D49").Value = '45bac4192c61cfeba14a64386a4b6e65'
$ws.Range("B50").Value = 'fio19.sh'
$ws.Range("D50").Value = '15f9d8fcfcd6ec32a51ab99a9e64f19e'
$ws.Range("B51").Value = 'fio20.sh'
$ws.Range("D51").Value = 'd3184ae79e15874650c264c7fc201803'
$ws.Range("B52").Value = 'fio18.sh'
$ws.Range("D52").Value = '2709ab002720bf259dfd7f79bf2a2677'
$ws.Range("B53").Value = 'fio21.sh'
$ws.Range("D53").Value = '619f1fbf4c45c242ba377bf6323dbc3b'
$ws.Range("B54").Value = 'fio13.sh'
$ws.Range("D54").Value = '020e043234bba88e36b5941af15254d5'
$ws.Range("B55").Value = 'fio12.sh'
$ws.Range("D55").Value = '5dd1be85dd4ab61a7dda604d81f0fc65'
$ws.Range("B56").Value = 'fio20.sh'
$ws.Range("D56").Value = 'd3184ae79e15874650c264c7fc201803'
$ws.Range("B57").Value = 'fio24.sh'
$ws.Range("D57").Value = '9889b12792dc5156d8b0bea1b6a9cc63'
$ws.Range("B58").Value = 'fio23.sh'
$ws.Range("D58").Value = 'de9caef7e5395d10cc8e56ee6c4e6d96'
$ws.Range("B59").Value = 'fio21.sh'
$ws.Range("D59").Value = '619f1fbf4c45c242ba377bf6323dbc3b'
$ws.Range("B60").Value = 'run14.txt'
$ws.Range("B61").Value = 'run8.txt'
$ws.Range("B62").Value = 'run10.txt'
$ws.Range("B63").Value = 'run7.txt'
$ws.Range("B64").Value = 'run5.txt'
$ws.Range("B65").Value = 'run20.txt'
$ws.Range("B66").Value = 'run15.txt'
$ws.Range("B67").Value = 'run22.txt'
$ws.Range("B68").Value = 'run19.txt'
$ws.Range("B69").Value = 'run2.txt'
$ws.Range("B70").Value = 'run6.txt'
$ws.Range("B71").Value = 'run21.txt'
$ws.Range("B73").Value = 'run23.txt'
$ws.Range("B74").Value = 'run18.txt'
$ws.Range("B75").Value = 'run17.txt'
$ws.Range("B76").Value = 'run11.txt'
$ws.Range("B77").Value = 'run13.txt'
$ws.Range("B78").Value = 'run4.txt'
$ws.Range("B79").Value = 'run3.txt'
$ws.Range("B80").Value = 'run16.txt'
$ws.Range("B81").Value = 'run9.txt'
$ws.Range("B82").Value = 'run12.txt'
$ws.Range("B83").Value = 'device.ios3'
$ws.Range("D83").Value = 'bf6d291ed6e2c0330b061c3ce2422e24'
$ws.Range("B84").Value = 'device.ios4'
$ws.Range("D84").Value = 'ad62b8f76d3623bf276e3206c5ce6424'
$ws.Range("B85").Value = 'device.ios2'
$ws.Range("D85").Value = '51b08f41de06d5bd0ff0259c5f0efc3f'
$ws.Range("B90").Value = 'devlist.vm3_4kall.bkp'
$ws.Range("D90").Value = '4fdde708a5376fc76f2391c330a9310c'
$ws.Range("E90").Value = 1
$ws.Range("B91").Value = 'devlist.vm4_4kall'
$ws.Range("D91").Value = 'c17a73aa1a4d93f49861b304eedbecae'
$ws.Range("E91").Value = 0
$ws.Range("B95").Value = 'KKD geo tagging.xlsx'
$ws.Range("B97").Value = 'KKD geo tagging (2).xlsx'
